# Day 7 PPT update: slide 10 -- remove the empty placeholder textbox that used
# to sit "in front of" the picture, resize/reposition the picture to fill the
# space, and add a small new (empty) caption textbox below-right of the image.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# EMU -> points helper. Nudge by half an EMU so that floating point
# point-values round-trip to the exact target EMU integer on save.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

# ---------------------------------------------------------------------------
# 1. Remove the empty text box (old id 171) that sat behind/around the image.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Google Shape;171;ga561df5dd9_1_12") {
        $shp.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Reposition / resize the picture (previously id 172) to be bigger.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Google Shape;172;ga561df5dd9_1_12") {
        $shp.Left = EmuToPt 3062638
        $shp.Top = EmuToPt 2336350
        $shp.Width = EmuToPt 6066775
        $shp.Height = EmuToPt 3856650
        $shp.Name = "Google Shape;171;ga561df5dd9_1_12"
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Add the new (empty) caption textbox after the picture.
# ---------------------------------------------------------------------------
$newLeft = EmuToPt 6046725
$newTop = EmuToPt 6169375
$newWidth = EmuToPt 811200
$newHeight = EmuToPt 400200

$newBox = $s.Shapes.AddTextbox(1, $newLeft, $newTop, $newWidth, $newHeight)
$newBox.Name = "Google Shape;172;ga561df5dd9_1_12"

$newBox.Fill.Visible = 0
$newBox.Line.Visible = 0

$tf = $newBox.TextFrame
$tf.WordWrap = -1
$tf.VerticalAnchor = 1
$tf.MarginLeft = EmuToPt 91425
$tf.MarginRight = EmuToPt 91425
$tf.MarginTop = EmuToPt 91425
$tf.MarginBottom = EmuToPt 91425
$tf.AutoSize = 1

# AutoSize recalculates the box height from the (empty) text -- put the
# geometry back to the exact target box afterwards.
$newBox.Left = $newLeft
$newBox.Top = $newTop
$newBox.Width = $newWidth
$newBox.Height = $newHeight

$tr = $tf.TextRange
$tr.ParagraphFormat.Alignment = 1
$tr.IndentLevel = 1
$tr.ParagraphFormat.SpaceBefore = 0
$tr.ParagraphFormat.SpaceAfter = 0
$tr.ParagraphFormat.Bullet.Visible = 0
$tr.Font.Name = "Gill Sans"
